$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text must be preserved exactly as text (numbers/percents) use a
# temporary text number-format so Excel does not coerce them to numeric values;
# the style is then reset to Normal so no extra style survives in the saved file.
function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue $ws "D2" "327.23"
Set-TextValue $ws "E2" "-1.78%"
Set-TextValue $ws "G2" "12"

# Row 3
Set-TextValue $ws "D3" "44.33"
Set-TextValue $ws "E3" "0.98%"
Set-TextValue $ws "G3" "12"

# Row 4
Set-TextValue $ws "D4" "5.569"
Set-TextValue $ws "E4" "-2.13%"
Set-TextValue $ws "G4" "12"

# Row 5
Set-TextValue $ws "D5" "0.08051"
Set-TextValue $ws "E5" "-3.63%"
Set-TextValue $ws "G5" "12"

# Row 6
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws "D6" "4.311"
Set-TextValue $ws "E6" "-4.76%"
Set-TextValue $ws "G6" "12"

# Row 7
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws "D7" "1.898"
Set-TextValue $ws "E7" "-3.04%"
Set-TextValue $ws "G7" "12"

# Row 8
Set-TextValue $ws "E8" "-7.58%"
Set-TextValue $ws "G8" "12"

# Row 9
Set-TextValue $ws "D9" "0.9470"
Set-TextValue $ws "E9" "-0.04%"
Set-TextValue $ws "G9" "12"

# Row 10
Set-TextValue $ws "D10" "0.1139"
Set-TextValue $ws "E10" "-7.27%"
Set-TextValue $ws "G10" "12"

# Row 11
Set-TextValue $ws "D11" "0.1835"
Set-TextValue $ws "E11" "-7.00%"
Set-TextValue $ws "G11" "12"

# Row 12
Set-TextValue $ws "D12" "0.09640"
Set-TextValue $ws "E12" "-3.99%"
Set-TextValue $ws "G12" "12"

# Row 13
Set-TextValue $ws "D13" "0.04369"
Set-TextValue $ws "E13" "-1.79%"
Set-TextValue $ws "G13" "12"

# Row 14
Set-TextValue $ws "D14" "0.1067"
Set-TextValue $ws "E14" "-0.12%"
Set-TextValue $ws "G14" "12"

# Row 15
Set-TextValue $ws "D15" "0.001283"
Set-TextValue $ws "E15" "-2.39%"
Set-TextValue $ws "G15" "12"

# Row 16
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D16" "0.005990"
Set-TextValue $ws "E16" "1.66%"
Set-TextValue $ws "G16" "12"

# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D17" "3.623"
Set-TextValue $ws "E17" "4.22%"
Set-TextValue $ws "G17" "12"

# Row 18
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue $ws "D18" "0.3498"
Set-TextValue $ws "E18" "-1.10%"
Set-TextValue $ws "G18" "12"

# Row 19
$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws "D19" "9.580"
Set-TextValue $ws "E19" "9.73%"
Set-TextValue $ws "G19" "12"

# Row 20
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue $ws "D20" "0.1379"
Set-TextValue $ws "E20" "1.22%"
Set-TextValue $ws "G20" "12"

# Row 21
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue $ws "D21" "0.2653"
Set-TextValue $ws "E21" "-2.38%"
Set-TextValue $ws "G21" "12"

# Row 22
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws "D22" "0.04222"
Set-TextValue $ws "E22" "-3.98%"
Set-TextValue $ws "G22" "12"

# Row 23
Set-TextValue $ws "E23" "0.93%"
Set-TextValue $ws "G23" "12"

# Row 24
Set-TextValue $ws "D24" "0.004469"
Set-TextValue $ws "E24" "2.66%"
Set-TextValue $ws "G24" "12"

# Row 25
Set-TextValue $ws "D25" "0.0001262"
Set-TextValue $ws "E25" "2.33%"
Set-TextValue $ws "G25" "12"

# Row 26
Set-TextValue $ws "D26" "0.0003994"
Set-TextValue $ws "E26" "-0.05%"
Set-TextValue $ws "G26" "12"

# Row 27
Set-TextValue $ws "G27" "12"

# Row 28
Set-TextValue $ws "G28" "12"

# Row 29
Set-TextValue $ws "G29" "12"

# Row 30
Set-TextValue $ws "G30" "12"

# Row 31
Set-TextValue $ws "G31" "12"

# Row 32
Set-TextValue $ws "G32" "12"

# Row 33
Set-TextValue $ws "G33" "12"

# Row 34
Set-TextValue $ws "G34" "12"

# Row 35
Set-TextValue $ws "G35" "12"

# Row 36
Set-TextValue $ws "G36" "12"

# Row 37
Set-TextValue $ws "G37" "12"

# Row 38
Set-TextValue $ws "D38" "0.02623"
Set-TextValue $ws "E38" "-8.77%"
Set-TextValue $ws "G38" "12"

# Row 39
Set-TextValue $ws "D39" "0.05454"
Set-TextValue $ws "E39" "-7.56%"
Set-TextValue $ws "G39" "12"

# Row 40
Set-TextValue $ws "D40" "0.007579"
Set-TextValue $ws "E40" "-4.24%"
Set-TextValue $ws "G40" "12"

# Row 41
Set-TextValue $ws "D41" "0.1394"
Set-TextValue $ws "E41" "-2.03%"
Set-TextValue $ws "G41" "12"

# Row 42
Set-TextValue $ws "D42" "0.007341"
Set-TextValue $ws "E42" "-19.18%"
Set-TextValue $ws "G42" "12"

# Row 43
Set-TextValue $ws "D43" "0.002018"
Set-TextValue $ws "E43" "-5.96%"
Set-TextValue $ws "G43" "12"

# Row 44
Set-TextValue $ws "D44" "0.008839"
Set-TextValue $ws "E44" "-10.16%"
Set-TextValue $ws "G44" "12"

# Row 45
Set-TextValue $ws "D45" "0.00006937"
Set-TextValue $ws "E45" "-9.20%"
Set-TextValue $ws "G45" "12"

# Row 46
Set-TextValue $ws "D46" "0.00000000751"
Set-TextValue $ws "E46" "-0.04%"
Set-TextValue $ws "G46" "12"

# Row 47
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextValue $ws "D47" "0.002273"
Set-TextValue $ws "E47" "-0.05%"
Set-TextValue $ws "G47" "12"

# Row 48
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue $ws "D48" "0.004189"
Set-TextValue $ws "E48" "31.41%"
Set-TextValue $ws "G48" "12"

# Row 49
Set-TextValue $ws "D49" "0.00002103"
Set-TextValue $ws "E49" "-0.04%"
Set-TextValue $ws "G49" "12"

# Row 50
Set-TextValue $ws "D50" "0.0002003"
Set-TextValue $ws "E50" "-0.04%"
Set-TextValue $ws "G50" "12"

# Row 51
Set-TextValue $ws "G51" "12"
